$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''43.482.27'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +1.39%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''2.331.22'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +1.69%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''  -0.34%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''312.71'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -0.51%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''108.14'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +4.42%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''0.628'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  +0.30%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '''  -0.22%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.618'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  +2.53%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''41.12'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  +4.62%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.0919'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  +1.26%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''8.53'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  +2.18%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = '''  -1.12%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = '''  +1.73%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''15.48'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  +1.75%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''2.686.25'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +1.79%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''2.325.16'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  +1.96%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''43.671.92'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +2.15%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''7.54'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +1.32%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = '''  +1.20%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''13.08'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -4.08%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''74.35'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  +0.87%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = '''  -3.12%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''267.78'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +1.03%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = '''  +3.49%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = '''  -0.22%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''7.65'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  +9.66%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''11.13'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  +2.95%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = '''  -1.52%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''39.81'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  +8.35%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''22.57'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +0.03%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''167.78'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  +0.44%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''0.0886'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  +1.77%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = '''  +8.54%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = '''  +0.24%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '''  +2.31%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''4.70'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  +3.57%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''0.0364'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +4.03%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''2.90'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +9.17%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''3.81'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  +3.48%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''1.71'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  +8.17%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''104.51'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  +10.03%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = 'Celestia'
$ws.Range("C43").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D43").Value = '''13.66'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  +13.83%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").Value = '''0.241'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  +4.47%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''71.67'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  +2.08%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''1.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +0.06%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''114.49'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +2.10%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = 'TheGraph'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D48").Value = '''0.220'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  +17.24%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '''1.660.08'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  -4.23%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = '''  +3.10%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''76.37'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  -4.84%  '
$ws.Range("E51").Style = "Normal"
